# "Generate Report for Handback" -- fills in the handback columns (Latest
# Target File / Latest Handback File / Latest Handback DateTime) for the
# zh-cn and de-de localization sheets now that the de-de handback has come
# back in sync with en-US, and widens a few columns so the new/longer
# content is readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$githubBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5a18e5c2f1e74136c35d6514b2bbf5ddcbc67b1b/e2e"

$statusHandedBack = "Handed back: in sync with en-US"
$zhcnHandbackTime = "2016-08-29 20:50:39"
$dedeHandbackTime = "2016-08-29 20:50:47"

# --- Status: flip from "Ready for handoff" to "Handed back: in sync with en-US" ---
$overview.Range("E2").Value = $statusHandedBack
$overview.Range("E3").Value = $statusHandedBack
$overview.Range("F2").Value = $statusHandedBack
$overview.Range("F3").Value = $statusHandedBack

$zhcn.Range("C2").Value = $statusHandedBack
$zhcn.Range("C3").Value = $statusHandedBack

$dede.Range("C2").Value = $statusHandedBack
$dede.Range("C3").Value = $statusHandedBack

# --- zh-cn: Latest Target File (I) / Latest Handback File (J) / Latest Handback DateTime (K) ---
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "$githubBase/a.md", "", "", "a.md")
$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K2").Value = $zhcnHandbackTime

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "$githubBase/a.md", "", "", "a.md")
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K3").Value = $zhcnHandbackTime

# --- de-de: Latest Target File (I) / Latest Handback File (J) / Latest Handback DateTime (K) ---
$dede.Hyperlinks.Add($dede.Range("I2"), "$githubBase/a.md", "", "", "a.md")
$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K2").Value = $dedeHandbackTime

$dede.Hyperlinks.Add($dede.Range("I3"), "$githubBase/a.md", "", "", "a.md")
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K3").Value = $dedeHandbackTime

# --- Widen columns that now hold longer handback status / handback-file text ---
# (character-width 29.15 / 39.1666667 round-trip to the same stored
# <col width=.../> the workbook ends up with after a manual column resize)
$overview.Columns.Item(5).ColumnWidth = 29.15   # E: zh-cn
$overview.Columns.Item(6).ColumnWidth = 29.15   # F: de-de

$zhcn.Columns.Item(3).ColumnWidth = 29.15        # C: Status
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664  # J: Latest Handback File

$dede.Columns.Item(3).ColumnWidth = 29.15        # C: Status
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664  # J: Latest Handback File

Write-Output "Generated handback report"
